# M14 Froze Token Embeddings
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Numeric updates in column C
$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 5
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 21
$ws.Range("C8").Value = 19
$ws.Range("C9").Value = 15
$ws.Range("C10").Value = 14
$ws.Range("C11").Value = 11
$ws.Range("C12").Value = 13
$ws.Range("C13").Value = 11
$ws.Range("C15").Value = 16
$ws.Range("C16").Value = 14
$ws.Range("C17").Value = 19
$ws.Range("C18").Value = 13

# Text updates in column B
$ws.Range("B10").Value = "<on>"
$ws.Range("B11").Value = "<make>"
$ws.Range("B16").Value = "<sile>"
